# Applies data-refresh updates to the 杭州-漫展信息 workbook
# (want-to-go counts and min ticket price columns), matching the
# "Update gh-pages to output generated at 456a3b4" commit.
$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 3080
$ws.Range("F6").Value = 2090
$ws.Range("F7").Value = 325
$ws.Range("F9").Value = 891
$ws.Range("F10").Value = 977
$ws.Range("F12").Value = 440
$ws.Range("F13").Value = 1144
$ws.Range("F16").Value = 533
$ws.Range("F17").Value = 7479
$ws.Range("F18").Value = 312
$ws.Range("F21").Value = 216
$ws.Range("F23").Value = 452
$ws.Range("F24").Value = 518
$ws.Range("F27").Value = 960
$ws.Range("F29").Value = 887
$ws.Range("F31").Value = 1135
$ws.Range("F32").Value = 1912
$ws.Range("F33").Value = 477
$ws.Range("F38").Value = 157
$ws.Range("F39").Value = 308
$ws.Range("F41").Value = 204

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value = 88

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("G3").Value = 88
$ws.Range("F7").Value = 3080
$ws.Range("F9").Value = 2090
$ws.Range("F10").Value = 325
$ws.Range("F12").Value = 891
$ws.Range("F14").Value = 977
$ws.Range("F16").Value = 440
$ws.Range("F17").Value = 1144
$ws.Range("F20").Value = 533
$ws.Range("F21").Value = 7479
$ws.Range("F22").Value = 312
$ws.Range("F26").Value = 216
$ws.Range("F28").Value = 452
$ws.Range("F29").Value = 518
$ws.Range("F32").Value = 960
$ws.Range("F34").Value = 887
$ws.Range("F36").Value = 1135
$ws.Range("F37").Value = 1912
$ws.Range("F38").Value = 477
$ws.Range("F43").Value = 157
$ws.Range("F44").Value = 308
$ws.Range("F49").Value = 204

